$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price cells to remain text (they are stored as
# inline/shared strings of numeric-looking text, e.g. "243.19") so that
# assigning the new value does not silently convert them to numbers.
$priceRange = $ws.Range("D2:D49")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "243.19"
$ws.Range("D3").Value = "23.04"
$ws.Range("D4").Value = "5.426"
$ws.Range("D5").Value = "0.05894"
$ws.Range("D6").Value = "3.449"
$ws.Range("D7").Value = "6.541"
$ws.Range("D8").Value = "0.8121"
$ws.Range("D9").Value = "0.9706"
$ws.Range("D10").Value = "0.1416"
$ws.Range("D11").Value = "0.07442"
$ws.Range("D12").Value = "0.03276"
$ws.Range("D13").Value = "0.03063"
$ws.Range("D14").Value = "0.09336"
$ws.Range("D15").Value = "3.868"
$ws.Range("D16").Value = "0.001575"
$ws.Range("D17").Value = "0.04678"
$ws.Range("D18").Value = "0.0005937"
$ws.Range("D19").Value = "0.005876"
$ws.Range("D20").Value = "0.001258"
$ws.Range("D24").Value = "2.135"
$ws.Range("D25").Value = "0.3230"
$ws.Range("D27").Value = "0.0002287"
$ws.Range("D40").Value = "0.03935"
$ws.Range("D41").Value = "0.006181"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D44").Value = "0.009137"
$ws.Range("D45").Value = "0.00005201"
$ws.Range("D47").Value = "0.7509"
$ws.Range("D48").Value = "0.002296"
$ws.Range("D49").Value = "0.00002103"

# Restore default (unformatted) cell appearance now that the values are set.
$priceRange.ClearFormats()
